$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 245.5
$ws.Range("I2").Value = 104.8
$ws.Range("K2").Value = 104.8
$ws.Range("M2").Value = 8.200000000000003

$ws.Range("H64").Value = 5047.8125
$ws.Range("I64").Value = 3942.7144
$ws.Range("J64").Value = 5907.3335
$ws.Range("K64").Value = 3942.7144
$ws.Range("L64").Value = 5907.3335
$ws.Range("M64").Value = -3694.7144
$ws.Range("N64").Value = -6403.3335

$ws.Range("H67").Value = 5047.8125
$ws.Range("I67").Value = 3942.7144
$ws.Range("J67").Value = 5907.3335
$ws.Range("K67").Value = 3942.7144
$ws.Range("L67").Value = 5907.3335
$ws.Range("M67").Value = -3084.7144
$ws.Range("N67").Value = -7623.3335

$ws.Range("H74").Value = 6750
$ws.Range("I74").Value = 3500
$ws.Range("K74").Value = 3500
$ws.Range("M74").Value = -2564

$ws.Range("H77").Value = 6750
$ws.Range("I77").Value = 3500
$ws.Range("K77").Value = 17500
$ws.Range("M77").Value = -12820

$ws.Range("H98").Value = 3236.7407
$ws.Range("I98").Value = 3458.5715
$ws.Range("K98").Value = 3458.5715
$ws.Range("M98").Value = -1960.5715

$ws.Range("H112").Value = 3862.4722
$ws.Range("I112").Value = 4067.4211
$ws.Range("J112").Value = 3633.4119
$ws.Range("K112").Value = 12202.2633
$ws.Range("L112").Value = 10900.2357
$ws.Range("M112").Value = -11094.2633
$ws.Range("N112").Value = -13116.2357

$ws.Range("H113").Value = 9960.571
$ws.Range("I113").Value = 13394.875
$ws.Range("J113").Value = 5381.5
$ws.Range("K113").Value = 13394.875
$ws.Range("L113").Value = 5381.5
$ws.Range("M113").Value = -10140.875
$ws.Range("N113").Value = -11889.5

$ws.Range("H116").Value = 5711.722
$ws.Range("J116").Value = 5163.1113
$ws.Range("L116").Value = 5163.1113
$ws.Range("N116").Value = -12047.1113

$ws.Range("H122").Value = 3236.7407
$ws.Range("I122").Value = 3458.5715
$ws.Range("K122").Value = 10375.7145
$ws.Range("M122").Value = -7925.7145

$ws.Range("H123").Value = 249999.33
$ws.Range("J123").Value = 249999.33
$ws.Range("L123").Value = 249999.33
$ws.Range("N123").Value = -259799.33

$ws.Range("H132").Value = 10909.667
$ws.Range("I132").Value = 7492.759
$ws.Range("K132").Value = 22478.277
$ws.Range("M132").Value = -19948.277

$ws.Range("H141").Value = 3337.6875
$ws.Range("I141").Value = 3273
$ws.Range("J141").Value = 3687
$ws.Range("K141").Value = 9819
$ws.Range("L141").Value = 11061
$ws.Range("M141").Value = -4639
$ws.Range("N141").Value = -21421

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 16508.334
$ws.Range("I34").Value = 9762.5
$ws.Range("J34").Value = 30000
$ws.Range("K34").Value = 9762.5
$ws.Range("L34").Value = 30000
$ws.Range("M34").Value = -9491.5
$ws.Range("N34").Value = -30542

$ws.Range("H39").Value = 1999.875
$ws.Range("I39").Value = 1714.1428
$ws.Range("J39").Value = 4000
$ws.Range("K39").Value = 1714.1428
$ws.Range("L39").Value = 4000
$ws.Range("M39").Value = -1194.1428
$ws.Range("N39").Value = -5040

$ws.Range("H45").Value = 1415.3
$ws.Range("I45").Value = 1242.5333
$ws.Range("J45").Value = 1933.6
$ws.Range("K45").Value = 1242.5333
$ws.Range("L45").Value = 1933.6
$ws.Range("M45").Value = -865.5333000000001
$ws.Range("N45").Value = -2687.6

$ws.Range("H61").Value = 2255976.5
$ws.Range("I61").Value = 6189.3
$ws.Range("K61").Value = 6189.3
$ws.Range("M61").Value = -5977.3

$ws.Range("H74").Value = 38002.855
$ws.Range("I74").Value = 4201.3
$ws.Range("J74").Value = 122506.75
$ws.Range("K74").Value = 4201.3
$ws.Range("L74").Value = 122506.75
$ws.Range("M74").Value = -3327.3
$ws.Range("N74").Value = -124254.75

$ws.Range("H77").Value = 38002.855
$ws.Range("I77").Value = 4201.3
$ws.Range("J77").Value = 122506.75
$ws.Range("K77").Value = 21006.5
$ws.Range("L77").Value = 612533.75
$ws.Range("M77").Value = -16638.5
$ws.Range("N77").Value = -621269.75

$ws.Range("H97").Value = 402.3846
$ws.Range("I97").Value = 444.5
$ws.Range("K97").Value = 444.5
$ws.Range("M97").Value = 51.5

$ws.Range("H132").Value = 4569851.5
$ws.Range("I132").Value = 2762.5
$ws.Range("K132").Value = 8287.5
$ws.Range("M132").Value = -5757.5

$ws.Range("H136").Value = 2255976.5
$ws.Range("I136").Value = 6189.3
$ws.Range("K136").Value = 18567.9
$ws.Range("M136").Value = -16017.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 30613.1
$ws.Range("I20").Value = 13338.308
$ws.Range("K20").Value = 13338.308
$ws.Range("M20").Value = -13091.308

$ws.Range("H80").Value = 2085.9
$ws.Range("I80").Value = 1076.4
$ws.Range("J80").Value = 2422.4
$ws.Range("K80").Value = 1076.4
$ws.Range("L80").Value = 2422.4
$ws.Range("M80").Value = -78.40000000000009
$ws.Range("N80").Value = -4418.4

$ws.Range("H83").Value = 2085.9
$ws.Range("I83").Value = 1076.4
$ws.Range("J83").Value = 2422.4
$ws.Range("K83").Value = 5382
$ws.Range("L83").Value = 12112
$ws.Range("M83").Value = -390
$ws.Range("N83").Value = -22096

$ws.Range("H94").Value = 2879
$ws.Range("J94").Value = 2008
$ws.Range("L94").Value = 2008
$ws.Range("N94").Value = -2910

$ws.Range("H134").Value = 151925.56
$ws.Range("I134").Value = 193883.83
$ws.Range("J134").Value = 68009
$ws.Range("K134").Value = 581651.49
$ws.Range("L134").Value = 204027
$ws.Range("M134").Value = -579116.49
$ws.Range("N134").Value = -209097

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 3268.5
$ws.Range("I32").Value = 3077.3333
$ws.Range("J32").Value = 3842
$ws.Range("K32").Value = 3077.3333
$ws.Range("L32").Value = 3842
$ws.Range("M32").Value = -2761.3333
$ws.Range("N32").Value = -4474

$ws.Range("H35").Value = 2260.625
$ws.Range("I35").Value = 1708.5526
$ws.Range("J35").Value = 12750
$ws.Range("K35").Value = 1708.5526
$ws.Range("L35").Value = 12750
$ws.Range("M35").Value = -1414.5526
$ws.Range("N35").Value = -13338

$ws.Range("H58").Value = 15352.619
$ws.Range("I58").Value = 6481.647
$ws.Range("J58").Value = 53054.25
$ws.Range("K58").Value = 6481.647
$ws.Range("L58").Value = 53054.25
$ws.Range("M58").Value = -6278.647
$ws.Range("N58").Value = -53460.25

$ws.Range("H132").Value = 81449630
$ws.Range("I132").Value = 2165.2
$ws.Range("K132").Value = 6495.599999999999
$ws.Range("M132").Value = -3965.599999999999

$ws.Range("H134").Value = 45461560
$ws.Range("I134").Value = 2247.2
$ws.Range("K134").Value = 6741.599999999999
$ws.Range("M134").Value = -4206.599999999999

$ws.Range("H136").Value = 15352.619
$ws.Range("I136").Value = 6481.647
$ws.Range("J136").Value = 53054.25
$ws.Range("K136").Value = 19444.941
$ws.Range("L136").Value = 159162.75
$ws.Range("M136").Value = -16894.941
$ws.Range("N136").Value = -164262.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 13133.333
$ws.Range("I110").Value = 2200
$ws.Range("K110").Value = 6600
$ws.Range("M110").Value = -2510

$ws.Range("H113").Value = 4587.357
$ws.Range("I113").Value = 10834
$ws.Range("J113").Value = 1117
$ws.Range("K113").Value = 32502
$ws.Range("L113").Value = 3351
$ws.Range("M113").Value = -30332
$ws.Range("N113").Value = -7691

$ws.Range("H138").Value = 2907.125
$ws.Range("I138").Value = 2204
$ws.Range("J138").Value = 5016.5
$ws.Range("K138").Value = 6612
$ws.Range("L138").Value = 15049.5
$ws.Range("M138").Value = -1472
$ws.Range("N138").Value = -25329.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 19166.666

$ws.Range("H97").Value = 1167.6
$ws.Range("I97").Value = 1035.0952
$ws.Range("K97").Value = 1035.0952
$ws.Range("M97").Value = -539.0952

$ws.Range("H105").Value = 21750
$ws.Range("J105").Value = 21750
$ws.Range("L105").Value = 21750
$ws.Range("N105").Value = -28738

$ws.Range("H132").Value = 1227839.6
$ws.Range("I132").Value = 11838.125
$ws.Range("K132").Value = 35514.375
$ws.Range("M132").Value = -32984.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 4912.2
$ws.Range("I32").Value = 3235.7778
$ws.Range("K32").Value = 3235.7778
$ws.Range("M32").Value = -2918.7778

$ws.Range("H40").Value = 6538.857
$ws.Range("I40").Value = 6193.5
$ws.Range("K40").Value = 6193.5
$ws.Range("M40").Value = -6057.5

$ws.Range("H82").Value = 2712.4546
$ws.Range("I82").Value = 2999.6924
$ws.Range("K82").Value = 2999.6924
$ws.Range("M82").Value = -2638.6924

$ws.Range("H85").Value = 2712.4546
$ws.Range("I85").Value = 2999.6924
$ws.Range("K85").Value = 2999.6924
$ws.Range("M85").Value = -1751.6924

$ws.Range("H132").Value = 1793337.8
$ws.Range("J132").Value = 4366696.5
$ws.Range("L132").Value = 13100089.5
$ws.Range("N132").Value = -13105149.5

$ws.Range("H136").Value = 1691029.6
$ws.Range("I136").Value = 27847
$ws.Range("J136").Value = 7927964.5
$ws.Range("K136").Value = 83541
$ws.Range("L136").Value = 23783893.5
$ws.Range("M136").Value = -80991
$ws.Range("N136").Value = -23788993.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 842233.9
$ws.Range("I132").Value = 4919.4287
$ws.Range("K132").Value = 14758.2861
$ws.Range("M132").Value = -12228.2861

$ws.Range("H136").Value = 402018.66
$ws.Range("I136").Value = 2236.1904
$ws.Range("K136").Value = 6708.5712
$ws.Range("M136").Value = -4158.5712

Write-Output "done"